$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "campus"
[void]$ws.Range("H1").Select()
